# Apply updated crypto price/volume data per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.117.43"
$ws.Range("E2").Value = "  -3.41%  "

$ws.Range("D3").Value = "'1.849.26"
$ws.Range("E3").Value = "  -2.29%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'0.7029"
$ws.Range("E5").Value = "  -4.92%  "

$ws.Range("D6").Value = "'237.81"
$ws.Range("E6").Value = "  -2.16%  "

$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("D8").Value = "'0.3034"
$ws.Range("E8").Value = "  -4.32%  "

$ws.Range("D9").Value = "'0.07514"
$ws.Range("E9").Value = "  +4.27%  "

$ws.Range("D10").Value = "'23.32"
$ws.Range("E10").Value = "  -6.49%  "

$ws.Range("D11").Value = "'0.08118"
$ws.Range("E11").Value = "  -2.94%  "

$ws.Range("D12").Value = "'0.7248"
$ws.Range("E12").Value = "  -4.71%  "

$ws.Range("D13").Value = "'1.840.63"
$ws.Range("E13").Value = "  -4.47%  "

$ws.Range("E14").Value = "  -4.25%  "

$ws.Range("D15").Value = "'88.89"
$ws.Range("E15").Value = "  -4.41%  "

$ws.Range("D16").Value = "'29.084.93"
$ws.Range("E16").Value = "  -3.58%  "

$ws.Range("D17").Value = "'5.769"
$ws.Range("E17").Value = "  -6.61%  "

$ws.Range("D18").Value = "'238.16"
$ws.Range("E18").Value = "  -4.93%  "

$ws.Range("D19").Value = "'13.05"
$ws.Range("E19").Value = "  -4.37%  "

$ws.Range("D20").Value = "'0.000007651"
$ws.Range("E20").Value = "  -2.81%  "

$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "'2.085.06"
$ws.Range("E22").Value = "  -3.86%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("D24").Value = "'7.536"
$ws.Range("E24").Value = "  -5.52%  "

$ws.Range("D25").Value = "'161.94"
$ws.Range("E25").Value = "  -1.44%  "

$ws.Range("D26").Value = "'8.972"
$ws.Range("E26").Value = "  -3.59%  "

$ws.Range("D27").Value = "'0.1455"
$ws.Range("E27").Value = "  -8.11%  "

$ws.Range("D28").Value = "'18.00"
$ws.Range("E28").Value = "  -4.04%  "

$ws.Range("D29").Value = "'1.941"
$ws.Range("E29").Value = "  -5.99%  "

$ws.Range("D30").Value = "'1.388"
$ws.Range("E30").Value = "  -6.31%  "

$ws.Range("D31").Value = "'4.529"
$ws.Range("E31").Value = "  -1.37%  "

$ws.Range("D32").Value = "'1.492"
$ws.Range("E32").Value = "  -2.93%  "

$ws.Range("D33").Value = "'3.986"
$ws.Range("E33").Value = "  -5.37%  "

$ws.Range("D34").Value = "'0.05141"
$ws.Range("E34").Value = "  -4.52%  "

$ws.Range("D35").Value = "'1.186"
$ws.Range("E35").Value = "  -5.49%  "

$ws.Range("D36").Value = "'1.035"
$ws.Range("E36").Value = "  +3.22%  "

$ws.Range("E37").Value = "  -9.54%  "

$ws.Range("D38").Value = "'2.660"
$ws.Range("E38").Value = "  -2.50%  "

$ws.Range("D39").Value = "'0.01874"
$ws.Range("E39").Value = "  -4.96%  "

$ws.Range("D40").Value = "'2.679"
$ws.Range("E40").Value = "  -3.09%  "

$ws.Range("D41").Value = "'0.9406"
$ws.Range("E41").Value = "  +8.16%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'1.082.91"
$ws.Range("E42").Value = "  -1.70%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.985"
$ws.Range("E43").Value = "  -1.54%  "

$ws.Range("D44").Value = "'0.4285"
$ws.Range("E44").Value = "  -6.12%  "

$ws.Range("D45").Value = "'69.84"
$ws.Range("E45").Value = "  -3.97%  "

$ws.Range("E46").Value = "  -0.29%  "

$ws.Range("D47").Value = "'102.22"
$ws.Range("E47").Value = "  -2.01%  "

$ws.Range("D48").Value = "'1.742"
$ws.Range("E48").Value = "  -6.80%  "

$ws.Range("D49").Value = "'1.975.39"
$ws.Range("E49").Value = "  -4.42%  "

$ws.Range("D50").Value = "'9.153"
$ws.Range("E50").Value = "  -4.50%  "

$ws.Range("D51").Value = "'7.038"
$ws.Range("E51").Value = "  -7.51%  "

